# culture_collection を MIxS (MIGS.ba.host-associated.4.0) から再度削除する
#
# The template's header row (row 15) lists one field name per column,
# each with an explanatory cell comment underneath it. Column U holds
# "culture_collection" (with its comment). This field is being removed
# again, so every column from V onward (values + comments) shifts one
# column to the left, and the now-unused trailing column (previously
# BW) is cleared out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = 15
$firstCol  = 21   # U  (culture_collection - the column being removed)
$lastCol   = 75   # BW (last used column before the edit)

function Get-ColLetter {
    param([int]$n)
    $letter = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $letter = [char](65 + $rem) + $letter
        $n = [int](($n - 1) / 26)
    }
    return $letter
}

function Set-CommentText {
    param($range, $text)
    $c = $range.Comment
    if ($c -eq $null) {
        $range.AddComment($text)
    } else {
        $c.Text($text)
    }
}

function Clear-CommentIfAny {
    param($range)
    $c = $range.Comment
    if ($c -ne $null) {
        $c.Delete()
    }
}

# Shift columns (firstCol .. lastCol-1) to take on the content (value +
# comment) that currently lives one column to the right, ascending so
# each source column is read before it gets overwritten.
for ($i = $firstCol; $i -lt $lastCol; $i++) {
    $destAddr = (Get-ColLetter $i) + $headerRow
    $srcAddr  = (Get-ColLetter ($i + 1)) + $headerRow

    $destRange = $ws.Range($destAddr)
    $srcRange  = $ws.Range($srcAddr)

    $destRange.Value = $srcRange.Value()

    $srcComment = $srcRange.Comment
    if ($srcComment -eq $null) {
        Clear-CommentIfAny $destRange
    } else {
        Set-CommentText $destRange ($srcComment.Text())
    }
}

# The old last column (BW15) is now a duplicate of BV15 - clear it out,
# content and comment alike, since the sheet is one column narrower.
$lastAddr = (Get-ColLetter $lastCol) + $headerRow
$lastRange = $ws.Range($lastAddr)
Clear-CommentIfAny $lastRange
$lastRange.ClearContents()
